$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.008212327957153
$ws.Range("B1").Value = 3.204925537109375
$ws.Range("C1").Value = 3.844389677047729
$ws.Range("D1").Value = 2.020413875579834
$ws.Range("E1").Value = 1.193008184432983

$wb.Save()
